# Append a new "2025-03-15" row (row 14) to every price sheet in the
# workbook, carrying forward the same price that was last recorded
# (row 13) on each sheet. Mirrors the existing Date/Price columns:
# both cells must stay plain text (like the rest of column A/B),
# not get auto-converted into an Excel date serial / number by COM.

$wb = $excel.ActiveWorkbook

# Sheet name -> price to repeat on the new row.
$updates = [ordered]@{
    "N-Dense"                  = "40"
    "N-Type"                   = "43"
    "N-type Wafer"              = "1.19"
    "Cell Topcon 183mm"        = "0.298"
    "Module Topcon 183mm"      = "0.1"
    "Silver Rear_side"         = "5,455"
    "Silver Busbar front-side" = "8,167"
    "Silver finger front-side" = "8,217"
    "USD_CNY"                  = "7.2637"
}

$newDate = "2025-03-15"
$newRow = 14

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $price = $updates[$sheetName]

    $dateCell = $ws.Cells.Item($newRow, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $newDate
    $dateCell.Style = "Normal"

    $priceCell = $ws.Cells.Item($newRow, 2)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = "Normal"
}
